$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employee record did not exist yet in row 7 -> append it
$ws.Range("B7").Value = "0112548988482"
$ws.Range("A7").Value = "yasmeen  ahmed"
$ws.Range("C7").Value = ". Net developer"

$ws.Range("D7").Value = "test5@ahliunited.com"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:test5@ahliunited.com") | Out-Null
# Hyperlinks.Add mints its own ad-hoc style; re-apply the workbook's existing
# "Hyperlink" formatted cell (D6) so D7 keeps the same style index as D2:D6.
$ws.Range("D6").Copy($ws.Range("D7"))
$ws.Range("D7").Value = "test5@ahliunited.com"

$ws.Range("E7").Value = "assuit"
$ws.Range("F7").Value = 2351
$ws.Range("G7").Formula = "=F7*1.25"
$ws.Range("H7").Value = "female"

# Leftover formatted-but-empty hyperlink-style cell below the new row
$ws.Range("D6").Copy($ws.Range("D8"))
$ws.Range("D8").ClearContents()

$ws.Range("C13").Select()
